$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 84) matching the existing table's pattern:
# A = date (stored as text, like the other rows), B = weekday, C/D = numbers.

# Temporarily force column A's cell to text format so the date-like string
# "2025/10/09" isn't auto-converted into a date serial number, then restore
# the cell's style to Normal so it matches the unstyled look of the other
# data rows (no explicit style index left behind).
$ws.Range("A84").NumberFormat = "@"
$ws.Range("A84").Value = "2025/10/09"
$ws.Range("A84").Style = "Normal"

$ws.Range("B84").Value = "木"
$ws.Range("C84").Value = 11
$ws.Range("D84").Value = 201
